$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 390; this shifts existing rows 390-476
# down to 391-477 and grows the used range from T476 to T477.
$ws.Rows(390).Insert()

# Populate the new row 390 with the new weekly price observation.
# Columns A,B,C,E,F,G,H,I,J,R are constant across this block, so copy them
# from the (now-shifted) neighboring row 391.
$ws.Range("A390").Value = 7
$ws.Range("B390").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C390").Value = "Ñuble"
$ws.Range("D390").Value = 45258
$ws.Range("E390").Value = 16
$ws.Range("F390").Value = "Fruta"
$ws.Range("G390").Value = 100103
$ws.Range("H390").Value = "Frutos de hueso (carozo)"
$ws.Range("I390").Value = 100103006
$ws.Range("J390").Value = "Nectarín"
$ws.Range("K390").Value = "Early Glo"
$ws.Range("L390").Value = "Primera"
$ws.Range("M390").Value = 100
$ws.Range("N390").Value = 14000
$ws.Range("O390").Value = 15000
$ws.Range("P390").Value = 14500
$ws.Range("Q390").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R390").Value = "Región de O'Higgins"
$ws.Range("S390").Value = 967
$ws.Range("T390").Value = 15
